$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.908.09"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "1.552.15"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("D5").Value = "'206.71"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'0.485"
$ws.Range("E6").Value = "  +0.38%  "
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "'21.71"
$ws.Range("E8").Value = "  +2.06%  "
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").Value = "'0.0585"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "1.772.95"
$ws.Range("E12").Value = "  +1.38%  "
$ws.Range("D13").Value = "1.557.07"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "26.883.92"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D18").Value = "'216.58"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").Value = "'9.19"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").Value = "'153.61"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").Value = "'6.61"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("D30").Value = "'0.0465"
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("E32").Value = "  -0.47%  "
$ws.Range("D33").Value = "1.418.88"
$ws.Range("E33").Value = "  +4.80%  "
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("D36").Value = "'0.960"
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("D38").Value = "'0.0165"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("D39").Value = "'0.524"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("D40").Value = "'0.805"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").Value = "'5.64"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").Value = "'2.26"
$ws.Range("E44").Value = "  +3.77%  "
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("D46").Value = "'1.75"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").Value = "1.687.39"
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("D48").Value = "'86.21"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0101"
$ws.Range("E50").Value = "  +4.11%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0958"
$ws.Range("E51").Value = "  +1.63%  "
